$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(4).ColumnWidth = 12.830729166666666
$ws.Columns.Item(5).ColumnWidth = 13.498697916666666
$ws.Columns.Item(9).ColumnWidth = 11.998697916666666
$ws.Columns.Item(20).ColumnWidth = 15.666666666666666
